# Kazakh_data.xlsx edit: normalize valency-class columns.
# For every row where column I ("X") still holds the placeholder "TR" and
# column J ("Y") is empty, set I to "NOM" and J to "ACC" (column K, the
# locus, already holds "TR" and is left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(5,9,10,16,17,19,20,21,27,28,29,30,32,34,37,38,40,41,42,44,45,47,50,51,56,58,61,64,67,68,69,70,71,72,73,76,77,79,86,87,89,92,93,94,97,101,102,103,104,106,107,108,109,110,111,116,120,122,125,127,129,130)

foreach ($r in $rows) {
    $ws.Range("I$r").Value = "NOM"
    $ws.Range("J$r").Value = "ACC"
}

# Match the author's final selection (cell M9 was left selected).
$ws.Range("M9").Select()
